$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 3
    3  = 1
    4  = 1
    5  = 0
    6  = 1
    7  = 1
    8  = 2
    9  = 0
    10 = 1
    11 = 2
    12 = 0
    13 = 1
    14 = 1
    15 = 1
    16 = 1
    17 = 1
    18 = 1
    19 = 1
    20 = 0
    21 = 1
    22 = 2
    23 = 1
    24 = 1
    25 = 0
    26 = 1
    27 = 1
    28 = 1
    29 = 0
    30 = 2
    31 = 2
    32 = 2
    33 = 1
    34 = 2
    35 = 0
    36 = 1
    37 = 1
    38 = 1
    39 = 1
    40 = 1
    41 = 1
    42 = 2
    43 = 0
    44 = 2
    45 = 1
    46 = 1
    47 = 0
    48 = 0
    49 = 3
    50 = 2
    51 = 1
    52 = 0
    53 = 0
    54 = 1
    55 = 2
    56 = 0
    57 = 1
    58 = 1
    59 = 1
    60 = 0
    61 = 1
    62 = 1
    63 = 1
    64 = 0
    65 = 1
    66 = 1
    67 = 2
    68 = 2
    69 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
